# Weekly fruit/vegetable update: a new week of Espárragos price records was
# added for "Vega Central Mapocho de Santiago". This shifts the existing
# data rows 45-65 down by three rows (to 48-68) and inserts three brand new
# rows (45-47) with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows before row 45, pushing old rows 45-65 down to 48-68.
$ws.Rows("45:47").Insert()

# Populate the 3 newly inserted rows with the new week's records.
$newRows = @(
    @{ Row = 45; D = 44488; H = "Sin especificar"; I = "Banquete"; J = 250; K = 1600; L = 1600; M = 1600; N = "`$/kilo"; O = "Provincia de Linares"; P = 1600; Q = 1 },
    @{ Row = 46; D = 44488; H = "Sin especificar"; I = "Primera";  J = 430; K = 1400; L = 1400; M = 1400; N = "`$/kilo"; O = "Provincia de Linares"; P = 1400; Q = 1 },
    @{ Row = 47; D = 44488; H = "Sin especificar"; I = "Segunda";  J = 160; K = 1200; L = 1200; M = 1200; N = "`$/kilo"; O = "Provincia de Linares"; P = 1200; Q = 1 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 9
    $ws.Cells.Item($row, 2).Value = "Vega Central Mapocho de Santiago"
    $ws.Cells.Item($row, 3).Value = "Metropolitana"
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = 13
    $ws.Cells.Item($row, 6).Value = 300000000
    $ws.Cells.Item($row, 7).Value = "Espárragos"
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = "Hortaliza"
}
